# adding new progress as of date 04 nov 2025
#
# For each training row (3-17) on the "Training Dashboard" sheet:
#   - column H ("PERIOD TO EXPIRE") drops by one day
#   - column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025
#
# Column I holds the date as plain text (not a real Excel date), so we
# write it through a text formula and immediately paste-special it back
# as a value; that avoids Excel's automatic "looks like a date" literal
# conversion (which would otherwise turn the cell into a date serial
# number and change its number format/style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$xlPasteValues = -4163

for ($row = 3; $row -le 17; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE

    $hCell.Value = $hCell.Value() - 1

    $iCell.Formula = "=""04-Nov-2025"""
    $iCell.Copy()
    $iCell.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = 0
